$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the data set. It belongs right
# before the current row 244, so insert a fresh row there which pushes the
# existing rows 244:259 down to 245:260.
$ws.Rows.Item(244).Insert()

# Populate the newly inserted row 244 with the new record's values.
$ws.Cells.Item(244, 1).Value = 7
$ws.Cells.Item(244, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(244, 3).Value = "Ñuble"
$ws.Cells.Item(244, 4).Value = 45013
$ws.Cells.Item(244, 5).Value = 16
$ws.Cells.Item(244, 6).Value = "Fruta"
$ws.Cells.Item(244, 7).Value = 100104
$ws.Cells.Item(244, 8).Value = "Frutos de pepita"
$ws.Cells.Item(244, 9).Value = 100104005
$ws.Cells.Item(244, 10).Value = "Pera"
$ws.Cells.Item(244, 11).Value = "Packham's Triumph"
$ws.Cells.Item(244, 12).Value = "Primera"
$ws.Cells.Item(244, 13).Value = 50
$ws.Cells.Item(244, 14).Value = 10000
$ws.Cells.Item(244, 15).Value = 10000
$ws.Cells.Item(244, 16).Value = 10000
$ws.Cells.Item(244, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(244, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(244, 19).Value = 556
$ws.Cells.Item(244, 20).Value = 18
